{"js": "// The edit:\n//   \"...realizar X tablas, cuyos datos...\"\n// becomes\n//   \"...realizar X tablas(detallando nombres), cuyos datos...\"\n// and the lone \"_GoBack\" bookmark (previously sitting in an empty paragraph\n// further down, right after the \"RELACIONES\" list item) is moved to sit\n// right after the newly inserted \"(detallando nombres)\" text (i.e. right\n// before the \", cuyos datos...\" that follows).\n\n// 1) Remove the old \"_GoBack\" bookmark from its current location. It marks\n//    an empty paragraph, so deleting its (zero-width) range removes just the\n//    bookmark markers and leaves the empty paragraph intact.\nconst oldBookmark = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\noldBookmark.load(\"isNullObject\");\nawait context.sync();\nif (!oldBookmark.isNullObject) {\n  oldBookmark.delete();\n}\n\n// 2) Insert \"(detallando nombres)\" right after the word \"tablas\" in the\n//    \"Con este objetivo...\" paragraph.\nconst body = context.document.body;\nconst found = body.search(\"tablas\", { matchCase: false, matchWholeWord: false });\nfound.load(\"text\");\nawait context.sync();\n\nconst target = found.items[0];\nconst insertedRange = target.insertText(\"(detallando nombres)\", \"After\");\nawait context.sync();\n\n// 3) Re-create the \"_GoBack\" bookmark immediately after the inserted text\n//    (i.e. right before the comma that starts \", cuyos datos...\").\nconst insertionPoint = insertedRange.getRange(\"End\");\ninsertionPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The edit:\n#   \"...realizar X tablas, cuyos datos...\"\n# becomes\n#   \"...realizar X tablas(detallando nombres), cuyos datos...\"\n# and the lone \"_GoBack\" bookmark (previously sitting in an empty paragraph\n# further down, right after the \"RELACIONES\" list item) is moved to sit\n# right after the newly inserted \"(detallando nombres)\" text (i.e. right\n# before the \", cuyos datos...\" that follows).\n\n$d = $word.ActiveDocument\n\n# 1) Remove the old \"_GoBack\" bookmark from its current location (it marks\n#    an otherwise-empty paragraph, so deleting it only removes the bookmark\n#    markers and leaves the empty paragraph in place).\n$existing = $d.Bookmarks.Item(\"_GoBack\")\n$existing.Delete()\n\n# 2) Find the word \"tablas\" in the \"Con este objetivo...\" paragraph and\n#    insert \"(detallando nombres)\" right after it.\n$rng = $d.Content\n$rng.Find.Execute(\"tablas\") | Out-Null\n$rng.Collapse(0)                               # wdCollapseEnd\n$rng.InsertAfter(\"(detallando nombres)\")\n\n# 3) Re-create the \"_GoBack\" bookmark immediately after the inserted text\n#    (i.e. right before the comma that starts \", cuyos datos...\").\n$rng.Collapse(0)                               # wdCollapseEnd, end of inserted text\n$d.Bookmarks.Add(\"_GoBack\", $rng)\n"}
